# The commit removes the (empty/unused) slide titled
# "Solutions pour l'oscillateur anharmonique" from the deck. Every other
# difference visible in the target OOXML (slide part renumbering, the
# rId shifts inside <p:sldIdLst>, and an incidental namespace tweak on
# an unrelated mc:AlternateContent block elsewhere in the deck) is just
# PowerPoint repacking the package on save - a natural side effect of
# deleting a slide, not a separate edit.

$p = $ppt.ActivePresentation

$targetTitle = "Solutions pour l’oscillateur anharmonique"
$targetIndex = -1

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $text = $shape.TextFrame.TextRange.Text
            if ($text -eq $targetTitle) {
                $targetIndex = $i
                break
            }
        }
    }
    if ($targetIndex -ne -1) {
        break
    }
}

if ($targetIndex -eq -1) {
    # Fallback: this slide is the 27th slide in the original deck.
    $targetIndex = 27
}

$p.Slides.Item($targetIndex).Delete()
